$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M7").Value = 1000.77
$ws1.Range("M14").Value = 11316.44
$ws1.Range("H27").Value = 221.4
$ws1.Range("I27").Value = 671.4299999999999
$ws1.Range("M55").Value = "23 de 53"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F7").Value = 2542.76
$ws2.Range("F14").Value = 12869.49
$ws2.Range("F27").Value = 5906.93
$ws2.Range("F55").Value = 98715.52

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D7").Value = 2009.7
$ws3.Range("E7").Value = 390.3
$ws3.Range("F7").Value = 0.837375

$ws3.Range("D8").Value = 1428.09
$ws3.Range("E8").Value = -428.0899999999999
$ws3.Range("F8").Value = 1.42809

$ws3.Range("D16").Value = 52777.95
$ws3.Range("E16").Value = 3281.75
$ws3.Range("F16").Value = 0.9414597295383315

$ws3.Range("D19").Value = 98715.51999999999
$ws3.Range("E19").Value = 18724.17064517915
$ws3.Range("F19").Value = 0.8405635220740615
